# Daily attendance processing - 2025-11-05 23:21:38
#
# Normalises the "Recorded By" (column G) entries on the active sheet:
# a few specific combinations of recorder names/emails had their order
# flipped (e.g. "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com").
# Only rows whose current text exactly matches one of the known
# before-values are touched; everything else (single recorders, or
# combinations already in the desired order) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exact old-value -> new-value replacements observed for column G
# ("Recorded By").
$replacements = @{
    "backup@backdoor.com, system, System" = "backup@backdoor.com, System, system";
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com";
    "dnasr281@gmail.com, admin@admin.com" = "admin@admin.com, dnasr281@gmail.com";
}

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

# Column G = "Recorded By"
$col = 7

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $col)
    $text = $cell.Text

    if ($text -eq $null) { continue }
    if (-not $replacements.ContainsKey($text)) { continue }

    $cell.Value = $replacements[$text]
}
